$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 3171.0908
$ws.Range("J53").Value = 10323.333
$ws.Range("L53").Value = 10323.333
$ws.Range("N53").Value = -11597.333

# Row 62
$ws.Range("H62").Value = 2850.3
$ws.Range("I62").Value = 1728.7778
$ws.Range("K62").Value = 1728.7778
$ws.Range("M62").Value = -1104.7778

# Row 65
$ws.Range("H65").Value = 2850.3
$ws.Range("I65").Value = 1728.7778
$ws.Range("K65").Value = 8643.889000000001
$ws.Range("M65").Value = -5523.889000000001

# Row 86
$ws.Range("H86").Value = 10626.909
$ws.Range("I86").Value = 1725
$ws.Range("J86").Value = 15713.714
$ws.Range("K86").Value = 1725
$ws.Range("L86").Value = 15713.714
$ws.Range("M86").Value = -602
$ws.Range("N86").Value = -17959.714

# Row 89
$ws.Range("H89").Value = 10626.909
$ws.Range("I89").Value = 1725
$ws.Range("J89").Value = 15713.714
$ws.Range("K89").Value = 8625
$ws.Range("L89").Value = 78568.57000000001
$ws.Range("M89").Value = -3009
$ws.Range("N89").Value = -89800.57000000001

# Row 125
$ws.Range("H125").Value = 489.15384
$ws.Range("I125").Value = 469
$ws.Range("J125").Value = 492.81818
$ws.Range("K125").Value = 4221
$ws.Range("L125").Value = 4435.36362
$ws.Range("M125").Value = -1761
$ws.Range("N125").Value = -9355.36362

# Row 127
$ws.Range("H127").Value = 943.4375
$ws.Range("I127").Value = 462.375
$ws.Range("J127").Value = 1424.5
$ws.Range("K127").Value = 1387.125
$ws.Range("L127").Value = 4273.5
$ws.Range("M127").Value = 3572.875
$ws.Range("N127").Value = -14193.5

# Row 132
$ws.Range("H132").Value = 3114
$ws.Range("I132").Value = 3114
$ws.Range("K132").Value = 9342
$ws.Range("M132").Value = -6812

# Row 137
$ws.Range("H137").Value = 86872.27
$ws.Range("I137").Value = 118855.65
$ws.Range("K137").Value = 356566.95
$ws.Range("M137").Value = -354016.95

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2858.65
$ws.Range("I61").Value = 2660.7297
$ws.Range("K61").Value = 2660.7297
$ws.Range("M61").Value = -2448.7297

# Row 74
$ws.Range("H74").Value = 24391610
$ws.Range("I74").Value = 30303646
$ws.Range("K74").Value = 30303646
$ws.Range("M74").Value = -30302772

# Row 77
$ws.Range("H77").Value = 24391610
$ws.Range("I77").Value = 30303646
$ws.Range("K77").Value = 151518230
$ws.Range("M77").Value = -151513862

# Row 132
$ws.Range("H132").Value = 11848.804
$ws.Range("I132").Value = 1686.5476
$ws.Range("J132").Value = 59272.668
$ws.Range("K132").Value = 5059.642800000001
$ws.Range("L132").Value = 177818.004
$ws.Range("M132").Value = -2529.642800000001
$ws.Range("N132").Value = -182878.004

# Row 136
$ws.Range("H136").Value = 2858.65
$ws.Range("I136").Value = 2660.7297
$ws.Range("K136").Value = 7982.1891
$ws.Range("M136").Value = -5432.1891

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3832.1714
$ws.Range("I134").Value = 3798.4119
$ws.Range("K134").Value = 11395.2357
$ws.Range("M134").Value = -8860.235700000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3273.8413
$ws.Range("I31").Value = 1501.6389
$ws.Range("J31").Value = 5636.778
$ws.Range("K31").Value = 1501.6389
$ws.Range("L31").Value = 5636.778
$ws.Range("M31").Value = -1206.6389
$ws.Range("N31").Value = -6226.778

# Row 34
$ws.Range("H34").Value = 3273.8413
$ws.Range("I34").Value = 1501.6389
$ws.Range("J34").Value = 5636.778
$ws.Range("K34").Value = 1501.6389
$ws.Range("L34").Value = 5636.778
$ws.Range("M34").Value = -1299.6389
$ws.Range("N34").Value = -6040.778

# Row 132
$ws.Range("H132").Value = 3237.9092
$ws.Range("I132").Value = 2090.111
$ws.Range("K132").Value = 6270.333
$ws.Range("M132").Value = -3740.333

# Row 141
$ws.Range("H141").Value = 32818.35
$ws.Range("J141").Value = 32818.35
$ws.Range("L141").Value = 32818.35
$ws.Range("N141").Value = -43178.35

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 752.23
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 752.23
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2256.69
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12336.69

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1995.3226
$ws.Range("I102").Value = 1779.2963
$ws.Range("K102").Value = 1779.2963
$ws.Range("M102").Value = -157.2963

# Row 103
$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

# Row 107
$ws.Range("H107").Value = 269.27777
$ws.Range("I107").Value = 257.33334
$ws.Range("K107").Value = 257.33334
$ws.Range("M107").Value = 1662.66666

# Row 113
$ws.Range("H113").Value = 12538.875
$ws.Range("I113").Value = 15735.167
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 15735.167
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = -13565.167
$ws.Range("N113").Value = -7290

# Row 132
$ws.Range("H132").Value = 11178.3
$ws.Range("I132").Value = 2986.3428
$ws.Range("J132").Value = 22647.04
$ws.Range("K132").Value = 8959.028399999999
$ws.Range("L132").Value = 67941.12
$ws.Range("M132").Value = -6429.028399999999
$ws.Range("N132").Value = -73001.12

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2871.2856
$ws.Range("I7").Value = 2627.0908
$ws.Range("K7").Value = 2627.0908
$ws.Range("M7").Value = -2515.0908

# Row 61
$ws.Range("H61").Value = 3777.111
$ws.Range("I61").Value = 1912.2
$ws.Range("K61").Value = 1912.2
$ws.Range("M61").Value = -1710.2

# Row 113
$ws.Range("H113").Value = 3777.111
$ws.Range("I113").Value = 1912.2
$ws.Range("K113").Value = 1912.2
$ws.Range("M113").Value = 257.8

# Row 126
$ws.Range("H126").Value = 2871.2856
$ws.Range("I126").Value = 2627.0908
$ws.Range("K126").Value = 7881.2724
$ws.Range("M126").Value = -5411.2724

# Row 132
$ws.Range("H132").Value = 209948.39
$ws.Range("I132").Value = 288479.88
$ws.Range("J132").Value = 3803.25
$ws.Range("K132").Value = 865439.64
$ws.Range("L132").Value = 11409.75
$ws.Range("M132").Value = -862909.64
$ws.Range("N132").Value = -16469.75

# Row 136
$ws.Range("H136").Value = 1988.5588
$ws.Range("I136").Value = 1834.75
$ws.Range("J136").Value = 4449.5
$ws.Range("K136").Value = 5504.25
$ws.Range("L136").Value = 13348.5
$ws.Range("M136").Value = -2954.25
$ws.Range("N136").Value = -18448.5

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 151583.33
$ws.Range("J15").Value = 151583.33
$ws.Range("L15").Value = 151583.33
$ws.Range("N15").Value = -152159.33

# Row 45
$ws.Range("H45").Value = 10000
$ws.Range("J45").Value = 10000
$ws.Range("L45").Value = 10000
$ws.Range("N45").Value = -10982

# Row 81
$ws.Range("H81").Value = 300
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 600
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 461
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 300
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 3000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 2304
$ws.Range("N84").ClearContents()

# Row 136
$ws.Range("H136").Value = 27167388
$ws.Range("J136").Value = 12860
$ws.Range("L136").Value = 38580
$ws.Range("N136").Value = -43680
